$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Selectors")

# Column H ("2 Found") mirrors the "y" marker already present in column G
# ("1 Found") for every data row (2-17), since selectors previously counted
# as "1 Found" are now also recognized as "2 Found" once extracted
# automatically from code blocks.
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 8).Value = "y"
}

# Update the selection to reflect the newly filled column H range.
$ws.Range("G2:G17").Select()
